$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.897.07"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "1.548.63"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.16"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.99"
$ws.Range("E8").Value = "  +1.51%  "

$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  -0.59%  "

$ws.Range("D12").Value = "1.769.80"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "1.550.89"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").Value = "26.903.36"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.62"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.06"
$ws.Range("E18").Value = "  +1.07%  "

$ws.Range("D19").Value = "0.0₃0704"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.40"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.94"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("E31").Value = "  -1.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  -0.51%  "

$ws.Range("E33").Value = "  +3.91%  "

$ws.Range("D34").Value = "1.409.32"
$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("E35").Value = "  +2.37%  "

$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.527"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("E42").Value = "  +3.82%  "

$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D47").Value = "1.684.00"
$ws.Range("E47").Value = "  -0.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.04"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("E49").Value = "  +1.34%  "

$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +4.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  +0.24%  "

